$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 12:22:00"
$wsZhCn.Range("E4").Value = "2016-03-13 12:22:00"
$wsZhCn.Range("H2").Value = "2016-03-13 12:22:34"
$wsZhCn.Range("H4").Value = "2016-03-13 12:22:34"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 12:22:04"
$wsDeDe.Range("E4").Value = "2016-03-13 12:22:04"
$wsDeDe.Range("H2").Value = "2016-03-13 12:22:41"
$wsDeDe.Range("H4").Value = "2016-03-13 12:22:41"
